$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update "Pre Experimental Phase" column (B) values per updated outcome measures
$ws.Range("B2").Value = 68
$ws.Range("B3").Value = 70
$ws.Range("B4").Value = 71
$ws.Range("B5").Value = 72
$ws.Range("B6").Value = 83
$ws.Range("B7").Value = 82
$ws.Range("B8").Value = 70
$ws.Range("B9").Value = 72
$ws.Range("B10").Value = 71
$ws.Range("B11").Value = 47
$ws.Range("B12").Value = 53
$ws.Range("B13").Value = 66
$ws.Range("B14").Value = 72
